$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 10.91740197455928
$ws.Range("C3").Value = 10.40125770870062
$ws.Range("C4").Value = 10.0742843395387
$ws.Range("C5").Value = 9.273852805345385
$ws.Range("C6").Value = 9.029233340261023
$ws.Range("C7").Value = 8.60712335139333
$ws.Range("C8").Value = 7.10714812568638
$ws.Range("C9").Value = 7.918199678812821
$ws.Range("C10").Value = 5.391191953244578
